$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.001.60"
$ws.Range("E2").Value = "  +3.14%  "
$ws.Range("D3").Value = "2.613.92"
$ws.Range("E3").Value = "  +1.53%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "571.92"
$ws.Range("E5").Value = "  +0.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.85"
$ws.Range("E6").Value = "  +0.32%  "
$ws.Range("E7").Value = "  -0.25%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.603"
$ws.Range("E8").Value = "  +1.16%  "
$ws.Range("D9").Value = "2.642.90"
$ws.Range("E9").Value = "  +2.43%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.55"
$ws.Range("E10").Value = "  -2.26%  "
$ws.Range("E11").Value = "  +3.36%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.155"
$ws.Range("E12").Value = "  -2.64%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.370"
$ws.Range("E13").Value = "  +6.95%  "
$ws.Range("D14").Value = "3.078.45"
$ws.Range("E14").Value = "  +1.91%  "
$ws.Range("D15").Value = "60.983.44"
$ws.Range("E15").Value = "  +3.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "23.51"
$ws.Range("E16").Value = "  +4.59%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000142"
$ws.Range("E17").Value = "  +3.30%  "
$ws.Range("D18").Value = "2.626.51"
$ws.Range("E18").Value = "  +1.87%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.33"
$ws.Range("E19").Value = "  +10.66%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.68"
$ws.Range("E20").Value = "  +3.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "349.43"
$ws.Range("E21").Value = "  +3.48%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.16"
$ws.Range("E22").Value = "  +14.54%  "
$ws.Range("E23").Value = "  +0.13%  "
$ws.Range("E24").Value = "  +14.23%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "64.16"
$ws.Range("E25").Value = "  -0.46%  "
$ws.Range("B26").Value = "Kaspa"
$ws.Range("C26").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.163"
$ws.Range("E26").Value = "  +1.32%  "
$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.993"
$ws.Range("E27").Value = "  -0.31%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.74"
$ws.Range("E28").Value = "  +6.56%  "
$ws.Range("D29").Value = "0.0₃0799"
$ws.Range("E29").Value = "  +2.17%  "
$ws.Range("E30").Value = "  +7.55%  "
$ws.Range("E31").Value = "  -0.15%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.34"
$ws.Range("E32").Value = "  +4.20%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "160.93"
$ws.Range("E33").Value = "  +1.48%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.55"
$ws.Range("E34").Value = "  +2.84%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.30"
$ws.Range("E35").Value = "  +6.30%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.965"
$ws.Range("E36").Value = "  +10.58%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.21"
$ws.Range("E37").Value = "  +4.83%  "
$ws.Range("E38").Value = "  +5.96%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "37.79"
$ws.Range("E39").Value = "  +1.59%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.859"
$ws.Range("E40").Value = "  -1.55%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.82"
$ws.Range("E41").Value = "  +3.64%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "299.43"
$ws.Range("E42").Value = "  +2.11%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "139.83"
$ws.Range("E43").Value = "  +9.79%  "
$ws.Range("E44").Value = "  +1.21%  "
$ws.Range("E45").Value = "  -0.35%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.607"
$ws.Range("E46").Value = "  +2.51%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0550"
$ws.Range("E47").Value = "  +2.62%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0242"
$ws.Range("E48").Value = "  +3.90%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.83"
$ws.Range("E49").Value = "  +7.08%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "10.70"
$ws.Range("E50").Value = "  +0.64%  "
$ws.Range("D51").Value = "2.051.73"
$ws.Range("E51").Value = "  +5.31%  "
